$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new header row above the existing data (shifts rows 1-4 to rows 2-5,
# and Excel automatically re-points the existing SUM formula's references).
$ws.Rows.Item(1).Insert()

# Header row (row 1). Values are written in the order that reproduces the
# target shared-string table ordering: "Stages" (0), "Times" (1), "Total" (2).
$ws.Range("B1").Value = "Stages"
$ws.Range("A1").Value = "Times"
$ws.Range("C1").Value = "Total"

# Append the four additional stage-duration rows (rows 6-9), using the same
# number format as the other stage rows (style index 2 / numFmtId 21).
$ws.Cells.Item(6, 1).NumberFormat = "h:mm:ss"
$ws.Cells.Item(6, 1).Value = 0.09336805555555555

$ws.Cells.Item(7, 1).NumberFormat = "h:mm:ss"
$ws.Cells.Item(7, 1).Value = 0.14040509259259260

$ws.Cells.Item(8, 1).NumberFormat = "h:mm:ss"
$ws.Cells.Item(8, 1).Value = 0.08457175925925926

$ws.Cells.Item(9, 1).NumberFormat = "h:mm:ss"
$ws.Cells.Item(9, 1).Value = 0.03965277777777778

# Re-point the total formula in B2 to cover all stage rows except the last,
# and add a running total in C2 covering every stage row.
$ws.Range("B2").NumberFormat = "[h]:mm:ss"
$ws.Range("B2").Formula = "=SUM(A2:A8)"

$ws.Range("C2").NumberFormat = "[h]:mm:ss"
$ws.Range("C2").Formula = "=SUM(A2:A9)"

# New helper total for the final stage alone.
$ws.Range("B3").NumberFormat = "[h]:mm:ss"
$ws.Range("B3").Formula = "=SUM(A9:A9)"

# Match the saved selection/active-cell state from the target workbook.
$ws.Range("J9").Select()

$wb.Save()
